$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
# (column B, "Status", rows 2 & 3 on both the zh-cn and de-de sheets)
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Shared (language independent) handoff/source hyperlink targets.
# ---------------------------------------------------------------------------
$mdUrl1  = "https://github.com/OpenLocalizationTest/oltest/blob/81f31bcb6dc0eef71a2360491a314e754cfdda11/e2e/24b49878-1168-41df-9209-d28d7c304a0f.md"
$mdUrl2  = "https://github.com/OpenLocalizationTest/oltest/blob/81f31bcb6dc0eef71a2360491a314e754cfdda11/e2e/f1119a0a-5a0e-4964-908d-19df625c50a4.md"
$cfgUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/81f31bcb6dc0eef71a2360491a314e754cfdda11/.localization-config"

$mdDisp1 = "24b49878-1168-41df-9209-d28d7c304a0f.md"
$mdDisp2 = "f1119a0a-5a0e-4964-908d-19df625c50a4.md"
$cfgDisp = ".localization-config"

# ===========================================================================
# Overview sheet - Status columns (B/C) mirror the same "Ready for handoff"
# text that lives on the language sheets, so they need the same update.
# ===========================================================================
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value2 = $newStatus
$wsOverview.Range("C2").Value2 = $newStatus
$wsOverview.Range("B3").Value2 = $newStatus
$wsOverview.Range("C3").Value2 = $newStatus

# ===========================================================================
# zh-cn sheet
# ===========================================================================
$ws = $wb.Worksheets.Item("zh-cn")

$xlfUrl1  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0d42b02092171424b0de25767e99edbdb0a319c1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/24b49878-1168-41df-9209-d28d7c304a0f.5d1d59fb21b8c08a1c21d2083a422304aba18e35.zh-cn.xlf"
$xlfUrl2  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0d42b02092171424b0de25767e99edbdb0a319c1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/f1119a0a-5a0e-4964-908d-19df625c50a4.df3a4e6a6af4098ef13c7aa20cb724752ea738b3.zh-cn.xlf"
$xlfDisp1 = "24b49878-1168-41df-9209-d28d7c304a0f.5d1d59fb21b8c08a1c21d2083a422304aba18e35.zh-cn.xlf"
$xlfDisp2 = "f1119a0a-5a0e-4964-908d-19df625c50a4.df3a4e6a6af4098ef13c7aa20cb724752ea738b3.zh-cn.xlf"
$handbackDateTime = "2016-01-18 02:03:31"

$ws.Range("B2").Value2 = $newStatus
$ws.Range("B3").Value2 = $newStatus

$ws.Range("E2").Value2 = $mdDisp1
$ws.Range("F2").Value2 = $xlfDisp1
$ws.Range("E3").Value2 = $mdDisp2
$ws.Range("F3").Value2 = $xlfDisp2

$ws.Range("G2").Value2 = $handbackDateTime
$ws.Range("G3").Value2 = $handbackDateTime

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $mdUrl1, "", "", $mdDisp1)
$ws.Hyperlinks.Add($ws.Range("C2"), $xlfUrl1, "", "", $xlfDisp1)
$ws.Hyperlinks.Add($ws.Range("E2"), $mdUrl1, "", "", $mdDisp1)
$ws.Hyperlinks.Add($ws.Range("F2"), $xlfUrl1, "", "", $xlfDisp1)
$ws.Hyperlinks.Add($ws.Range("A3"), $mdUrl2, "", "", $mdDisp2)
$ws.Hyperlinks.Add($ws.Range("C3"), $xlfUrl2, "", "", $xlfDisp2)
$ws.Hyperlinks.Add($ws.Range("E3"), $mdUrl2, "", "", $mdDisp2)
$ws.Hyperlinks.Add($ws.Range("F3"), $xlfUrl2, "", "", $xlfDisp2)
$ws.Hyperlinks.Add($ws.Range("A4"), $cfgUrl, "", "", $cfgDisp)

# ===========================================================================
# de-de sheet
# ===========================================================================
$ws = $wb.Worksheets.Item("de-de")

$xlfUrl1  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a9332eb2c61478411e65d00ffa2f197d0788f0d4/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/24b49878-1168-41df-9209-d28d7c304a0f.5d1d59fb21b8c08a1c21d2083a422304aba18e35.de-de.xlf"
$xlfUrl2  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a9332eb2c61478411e65d00ffa2f197d0788f0d4/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/f1119a0a-5a0e-4964-908d-19df625c50a4.df3a4e6a6af4098ef13c7aa20cb724752ea738b3.de-de.xlf"
$xlfDisp1 = "24b49878-1168-41df-9209-d28d7c304a0f.5d1d59fb21b8c08a1c21d2083a422304aba18e35.de-de.xlf"
$xlfDisp2 = "f1119a0a-5a0e-4964-908d-19df625c50a4.df3a4e6a6af4098ef13c7aa20cb724752ea738b3.de-de.xlf"
$handbackDateTime = "2016-01-18 02:03:52"

$ws.Range("B2").Value2 = $newStatus
$ws.Range("B3").Value2 = $newStatus

$ws.Range("E2").Value2 = $mdDisp1
$ws.Range("F2").Value2 = $xlfDisp1
$ws.Range("E3").Value2 = $mdDisp2
$ws.Range("F3").Value2 = $xlfDisp2

$ws.Range("G2").Value2 = $handbackDateTime
$ws.Range("G3").Value2 = $handbackDateTime

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $mdUrl1, "", "", $mdDisp1)
$ws.Hyperlinks.Add($ws.Range("C2"), $xlfUrl1, "", "", $xlfDisp1)
$ws.Hyperlinks.Add($ws.Range("E2"), $mdUrl1, "", "", $mdDisp1)
$ws.Hyperlinks.Add($ws.Range("F2"), $xlfUrl1, "", "", $xlfDisp1)
$ws.Hyperlinks.Add($ws.Range("A3"), $mdUrl2, "", "", $mdDisp2)
$ws.Hyperlinks.Add($ws.Range("C3"), $xlfUrl2, "", "", $xlfDisp2)
$ws.Hyperlinks.Add($ws.Range("E3"), $mdUrl2, "", "", $mdDisp2)
$ws.Hyperlinks.Add($ws.Range("F3"), $xlfUrl2, "", "", $xlfDisp2)
$ws.Hyperlinks.Add($ws.Range("A4"), $cfgUrl, "", "", $cfgDisp)

Write-Host "Handback report generated."
